$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "244.07"
Set-TextValue "D3" "23.80"
Set-TextValue "D4" "5.254"
Set-TextValue "D5" "0.05834"
Set-TextValue "D6" "6.455"
Set-TextValue "D7" "3.332"
Set-TextValue "D8" "0.8080"
Set-TextValue "D9" "0.8958"
Set-TextValue "D10" "0.1381"
Set-TextValue "D11" "0.07099"
Set-TextValue "D12" "0.03064"
Set-TextValue "D13" "0.03026"
Set-TextValue "D14" "0.09325"
Set-TextValue "D15" "3.819"
Set-TextValue "D16" "0.001531"

Set-TextValue "D18" "0.0006048"
$ws.Range("E18").Value = "17OneONE"

Set-TextValue "D19" "0.006170"
Set-TextValue "D20" "0.001258"
Set-TextValue "D21" "0.004064"

Set-TextValue "D24" "2.175"

Set-TextValue "D40" "0.03815"
Set-TextValue "D41" "0.006268"
Set-TextValue "D42" "0.1053"
Set-TextValue "D43" "0.002531"
Set-TextValue "D44" "0.007001"
Set-TextValue "D45" "0.00005317"

Set-TextValue "D47" "0.4851"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

Set-TextValue "D48" "0.006318"
